# Daily attendance processing - 2025-12-12 08:36:12
#
# Column G ("Recorded By") lists the accounts that recorded/edited each
# attendance session, comma-separated. Historically "System" was always
# written first; this pass moves "System" to the end of the list (so the
# real/human recorder reads first) for every row where it currently leads.
# Rows whose value is just "System" alone, or that don't start with
# "System, " at all, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "System, "

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value()

    if ($value -ne $null -and $value.StartsWith($prefix)) {
        $remainder = $value.Substring($prefix.Length)
        $cell.Value = $remainder + ", System"
    }
}
